$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G14: the "slightly improves..." note gets an UPDATE appended ---
$ws.Range("G14").Value = "slightly improves over random forest | UPDATE: same as above"

# --- G13: the "significantly improves..." note gets an UPDATE appended ---
$ws.Range("G13").Value = "significantly improves over polynomial regression | UPDATE: it doesn’t actually improve anything, tree based methods are not for time series predictions, i.e they don’t perform well oustide training boundaries; Try something else."

# --- Row 13: new note "cancelled (see update)" in B13 ---
$ws.Range("B13").Value = "cancelled (see update)"

# --- Row 14: new note "do" in B14 (mirrors the "do" ditto-marks used elsewhere in the sheet) ---
$ws.Range("B14").Value = "do"

# --- Reset the view: clear the scrolled-to-H1 position and select B18 instead of L5 ---
$ws.Range("A1").Select()
$ws.Range("B18").Select()
